# Fix order of level in potential outcomes in the tot-tut table.
#
# Row 10 (labelled E[Y1]) and row 12 (labelled E[Y0]) each pull their four
# numbers (columns B:E) from the external "tot_tut" source workbook via
# formulas such as "=[1]tot_tut!B14" / "=[1]tot_tut!B17". The two source
# rows (14 and 17) had their treatment levels in the wrong order, so the
# numbers that were showing up under E[Y1] actually belonged to E[Y0] and
# vice-versa. The fix swaps which source row each displayed row pulls
# from, so the correct potential-outcome values line up with their label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ($\mathbb{E}[Y_1]$) used to read row 14 of the source table;
# it should read row 17 instead.
$ws.Range("B10").Formula = "=[1]tot_tut!B17"
$ws.Range("C10").Formula = "=[1]tot_tut!C17"
$ws.Range("D10").Formula = "=[1]tot_tut!D17"
$ws.Range("E10").Formula = "=[1]tot_tut!E17"

# Row 12 ($\mathbb{E}[Y_0]$) used to read row 17 of the source table;
# it should read row 14 instead.
$ws.Range("B12").Formula = "=[1]tot_tut!B14"
$ws.Range("C12").Formula = "=[1]tot_tut!C14"
$ws.Range("D12").Formula = "=[1]tot_tut!D14"
$ws.Range("E12").Formula = "=[1]tot_tut!E14"
